# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-price / profit data refresh described in the commit
# (updates currentAveragePrice* / LevePrice* / LeveProfit* columns on each class sheet).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 33.666668
$ws.Range("I31").Value = 33.666668
$ws.Range("K31").Value = 101.000004
$ws.Range("M31").Value = 128.999996
$ws.Range("H52").Value = 2172.6667
$ws.Range("I52").Value = 759
$ws.Range("K52").Value = 2277
$ws.Range("M52").Value = -2117
$ws.Range("H53").Value = 135.92308
$ws.Range("I53").Value = 57.625
$ws.Range("J53").Value = 261.2
$ws.Range("K53").Value = 57.625
$ws.Range("L53").Value = 261.2
$ws.Range("M53").Value = 579.375
$ws.Range("N53").Value = -1535.2
$ws.Range("H115").Value = 4190
$ws.Range("I115").Value = 4190
$ws.Range("K115").Value = 12570
$ws.Range("M115").Value = -11003
$ws.Range("H116").Value = 4259.8
$ws.Range("I116").Value = 3799.5
$ws.Range("J116").Value = 4566.6665
$ws.Range("K116").Value = 3799.5
$ws.Range("L116").Value = 4566.6665
$ws.Range("M116").Value = -357.5
$ws.Range("N116").Value = -11450.6665
$ws.Range("H137").Value = 2603.6843
$ws.Range("I137").Value = 1884.4667
$ws.Range("J137").Value = 5300.75
$ws.Range("K137").Value = 5653.4001
$ws.Range("L137").Value = 15902.25
$ws.Range("M137").Value = -3103.4001
$ws.Range("N137").Value = -21002.25
$ws.Range("H141").Value = 4232.778
$ws.Range("I141").Value = 3212.8572
$ws.Range("K141").Value = 9638.571599999999
$ws.Range("M141").Value = -4458.571599999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1375.4546
$ws.Range("I2").Value = 1375.4546
$ws.Range("K2").Value = 1375.4546
$ws.Range("M2").Value = -1262.4546
$ws.Range("H32").Value = 8562.237999999999
$ws.Range("I32").Value = 5843.606
$ws.Range("K32").Value = 5843.606
$ws.Range("M32").Value = -5556.606
$ws.Range("H45").Value = 1995.6666
$ws.Range("J45").Value = 2043.5
$ws.Range("L45").Value = 2043.5
$ws.Range("N45").Value = -2797.5
$ws.Range("H74").Value = 2007.6111
$ws.Range("I74").Value = 1186.6
$ws.Range("K74").Value = 1186.6
$ws.Range("M74").Value = -312.5999999999999
$ws.Range("H77").Value = 2007.6111
$ws.Range("I77").Value = 1186.6
$ws.Range("K77").Value = 5933
$ws.Range("M77").Value = -1565
$ws.Range("H116").Value = 1375.4546
$ws.Range("I116").Value = 1375.4546
$ws.Range("K116").Value = 1375.4546
$ws.Range("M116").Value = 918.5454
$ws.Range("H132").Value = 1830.0883
$ws.Range("I132").Value = 1766.8572
$ws.Range("K132").Value = 5300.571599999999
$ws.Range("M132").Value = -2770.571599999999
$ws.Range("H140").Value = 80000
$ws.Range("I140").Value = 65000
$ws.Range("J140").Value = 87500
$ws.Range("K140").Value = 65000
$ws.Range("L140").Value = 87500
$ws.Range("M140").Value = -59820
$ws.Range("N140").Value = -97860

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1375.4546
$ws.Range("I3").Value = 1375.4546
$ws.Range("K3").Value = 1375.4546
$ws.Range("M3").Value = -1261.4546
$ws.Range("H105").Value = 4462.5557
$ws.Range("J105").Value = 1675
$ws.Range("L105").Value = 1675
$ws.Range("N105").Value = -5169

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 538.2
$ws.Range("I22").Value = 599
$ws.Range("J22").Value = 497.66666
$ws.Range("K22").Value = 599
$ws.Range("L22").Value = 497.66666
$ws.Range("M22").Value = -249
$ws.Range("N22").Value = -1197.66666
$ws.Range("H31").Value = 5641.1055
$ws.Range("I31").Value = 4619.6665
$ws.Range("J31").Value = 7392.143
$ws.Range("K31").Value = 4619.6665
$ws.Range("L31").Value = 7392.143
$ws.Range("M31").Value = -4324.6665
$ws.Range("N31").Value = -7982.143
$ws.Range("H34").Value = 5641.1055
$ws.Range("I34").Value = 4619.6665
$ws.Range("J34").Value = 7392.143
$ws.Range("K34").Value = 4619.6665
$ws.Range("L34").Value = 7392.143
$ws.Range("M34").Value = -4417.6665
$ws.Range("N34").Value = -7796.143
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H58").Value = 3278.25
$ws.Range("I58").Value = 1275.909
$ws.Range("J58").Value = 4972.5386
$ws.Range("K58").Value = 1275.909
$ws.Range("L58").Value = 4972.5386
$ws.Range("M58").Value = -1072.909
$ws.Range("N58").Value = -5378.5386
$ws.Range("H86").Value = 10577.4
$ws.Range("I86").Value = 8971.75
$ws.Range("K86").Value = 8971.75
$ws.Range("M86").Value = -7848.75
$ws.Range("H89").Value = 10577.4
$ws.Range("I89").Value = 8971.75
$ws.Range("K89").Value = 44858.75
$ws.Range("M89").Value = -39242.75
$ws.Range("H94").Value = 365
$ws.Range("J94").Value = 221.33333
$ws.Range("L94").Value = 221.33333
$ws.Range("N94").Value = -1123.33333
$ws.Range("H107").Value = 803.3333
$ws.Range("I107").Value = 803.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 803.3333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1116.6667
$ws.Range("N107").ClearContents()
$ws.Range("H122").Value = 2998.5
$ws.Range("I122").Value = 3621.2307
$ws.Range("J122").Value = 1842
$ws.Range("K122").Value = 10863.6921
$ws.Range("L122").Value = 5526
$ws.Range("M122").Value = -8413.6921
$ws.Range("N122").Value = -10426
$ws.Range("H132").Value = 2364.6
$ws.Range("I132").Value = 2140.0588
$ws.Range("K132").Value = 6420.176399999999
$ws.Range("M132").Value = -3890.176399999999
$ws.Range("H136").Value = 3278.25
$ws.Range("I136").Value = 1275.909
$ws.Range("J136").Value = 4972.5386
$ws.Range("K136").Value = 3827.727
$ws.Range("L136").Value = 14917.6158
$ws.Range("M136").Value = -1277.727
$ws.Range("N136").Value = -20017.6158
$ws.Range("H141").Value = 132294.14
$ws.Range("I141").Value = 88798
$ws.Range("J141").Value = 149692.6
$ws.Range("K141").Value = 88798
$ws.Range("L141").Value = 149692.6
$ws.Range("M141").Value = -83618
$ws.Range("N141").Value = -160052.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25088656
$ws.Range("I4").Value = 29834844
$ws.Range("J4").Value = 1659.5714
$ws.Range("K4").Value = 89504532
$ws.Range("L4").Value = 4978.7142
$ws.Range("M4").Value = -89504420
$ws.Range("N4").Value = -5202.7142

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H99").Value = 12212.2
$ws.Range("I99").Value = 2187
$ws.Range("K99").Value = 2187
$ws.Range("M99").Value = 59
$ws.Range("H107").Value = 490.2
$ws.Range("I107").Value = 490.2
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 490.2
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1429.8
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 4375
$ws.Range("I113").Value = 4375
$ws.Range("K113").Value = 4375
$ws.Range("M113").Value = -2205
$ws.Range("H118").Value = 11000
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 3557.5
$ws.Range("I132").Value = 2502.625
$ws.Range("K132").Value = 7507.875
$ws.Range("M132").Value = -4977.875

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2777.8
$ws.Range("I68").Value = 2777.8
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2777.8
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2028.8
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2777.8
$ws.Range("I71").Value = 2777.8
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13889
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -10145
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 1203.6666
$ws.Range("I100").Value = 847.1667
$ws.Range("K100").Value = 847.1667
$ws.Range("M100").Value = -306.1667
$ws.Range("H116").Value = 106393.336
$ws.Range("J116").Value = 106393.336
$ws.Range("L116").Value = 106393.336
$ws.Range("N116").Value = -115571.336
$ws.Range("H117").Value = 35000
$ws.Range("J117").Value = 35000
$ws.Range("L117").Value = 35000
$ws.Range("N117").Value = -44178
$ws.Range("H132").Value = 3703.4055
$ws.Range("I132").Value = 2910.5
$ws.Range("K132").Value = 8731.5
$ws.Range("M132").Value = -6201.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 50000
$ws.Range("I41").Value = 50000
$ws.Range("K41").Value = 50000
$ws.Range("M41").Value = -49610
$ws.Range("H81").Value = 4775.2
$ws.Range("J81").Value = 789.5
$ws.Range("L81").Value = 1579
$ws.Range("N81").Value = -3701
$ws.Range("H84").Value = 4775.2
$ws.Range("J84").Value = 789.5
$ws.Range("L84").Value = 7895
$ws.Range("N84").Value = -18503
$ws.Range("H132").Value = 1591.5
$ws.Range("I132").Value = 1591.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4774.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2244.5
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 1463.1923
$ws.Range("I136").Value = 1463.1923
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4389.5769
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1839.5769
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 30000
$ws.Range("J139").Value = 0
$ws.Range("N139").ClearContents()

Write-Output "Applied scheduled Sheets update."